# Fruta / hortaliza, semanal
# Weekly update: a new daily price record is inserted at the top of the
# data block (row 162), pushing all the existing Pomelo records down by
# one row (162-240 -> 163-241).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 162; this shifts rows 162:240 down to 163:241 and
# grows the sheet dimension to A1:T241 automatically.
$ws.Rows("162:162").Insert()

# Populate the newly inserted row with this week's record.
$ws.Range("A162").Value = 10
$ws.Range("B162").Value = "Vega Modelo de Temuco"
$ws.Range("C162").Value = "La Araucanía"
$ws.Range("D162").Value = 44726
$ws.Range("E162").Value = 9
$ws.Range("F162").Value = "Fruta"
$ws.Range("G162").Value = 100102
$ws.Range("H162").Value = "Cítricos"
$ws.Range("I162").Value = 100102006
$ws.Range("J162").Value = "Pomelo"
$ws.Range("K162").Value = "Start Ruby"
$ws.Range("L162").Value = "Primera"
$ws.Range("M162").Value = 100
$ws.Range("N162").Value = 15000
$ws.Range("O162").Value = 15000
$ws.Range("P162").Value = 15000
$ws.Range("Q162").Value = "$/bandeja 15 kilos granel"
$ws.Range("R162").Value = "Región de O'Higgins"
$ws.Range("S162").Value = 1000
$ws.Range("T162").Value = 15
